$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- 1. Font rename: TimesNewToman -> Times New Roman (whole document) ---
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

# --- 2. Title / author / email ---
Replace-Text "Quantum Paradox: Unveiling the Enigma" "History: A Tapestry of Time Unraveled"
Replace-Text "Alex Mendez" "Anya Martinez"
Replace-Text "alex" "anya"
Replace-Text "mendez@emailhub" "martinez@eduworldmail"

# --- 3. First body paragraph (essay) ---
Replace-Text "In the realm of modern physics, the concept of quantum mechanics has opened doors to a captivating world of paradoxical phenomena that challenge our conventional understanding of reality" "As we venture through the annals of time, we find ourselves captivated by history's rich tapestry, intricately woven with the threads of human endeavor"
Replace-Text " Quantum entanglement, superposition, and the observer effect are just a few of these perplexities that continue to baffle physicists and philosophers alike" " Its vast canvas holds the vibrant hues of civilizations, empires, and cultures, each leaving an imprint on the ever-evolving story of our world"
Replace-Text " This essay delves into the enigmatic nature of quantum theory, exploring the paradoxes that have profoundly reshaped our comprehension of the universe" " From the grandeur of ancient empires to the complexities of modern societies, history offers us a profound window into the human experience, shaping our understanding of ourselves and our place in the grand scheme of things"
Replace-Text "Unraveling the complexities of entanglement, we encounter particles exhibiting a profound interconnectedness, sharing properties and instantaneous communication across vast distances" "Like a master storyteller, history weaves a narrative of cause and effect, revealing the intricate connections between events and their far-reaching consequences"
Replace-Text " The Schrodinger's cat paradox confronts us with the quandary of a cat's simultaneous existence in both alive and dead states within a sealed box until an observer opens it" " It unravels the tales of individuals whose actions reverberated through time, shaping destinies and altering the course of nations"
Replace-Text " The uncertainties inherent in quantum systems, exemplified by Heisenberg's uncertainty principle, blur the boundaries between particles and waves, challenging our notions of determinism and predictability" " Their trials, triumphs, and sacrifices become signposts along history's path, reminding us of the enduring impact of human agency"
Replace-Text "As we delve deeper into this enigmatic realm, we encounter the perplexing observer effect" "History also serves as a mirror to our present, reflecting our hopes, fears, and aspirations"
Replace-Text " This phenomenon suggests that the act of observation itself influences the behavior of subatomic particles, blurring the distinction between the observed and the observer" " By delving into the past, we gain a deeper understanding of the challenges and opportunities we face today"
Replace-Text " The paradoxical nature of quantum theory raises fundamental questions about the nature of reality, the limits of our knowledge, and the interconnectedness of the universe, stimulating ongoing debate and exploration" " It's a dialogue between the past and the present, where lessons learned and mistakes made offer invaluable insights for navigating the complexities of our own time"

# --- 4. Summary paragraph ---
Replace-Text "Quantum mechanics introduces a fascinating realm of paradoxes that challenge our classical intuitions and redefine our understanding of reality" "History is an immersive journey through the tapestry of time, offering a profound understanding of the human experience"
Replace-Text " Phenomena like entanglement, superposition, and the observer effect blur the lines between particles and waves, defy notions of determinism, and question the role of observers in shaping the outcomes of quantum systems" " It chronicles the rise and fall of civilizations, the struggles and triumphs of individuals, and the intricate connections between cause and effect"
Replace-Text " These paradoxes have spurred profound philosophical discussions, inspiring new perspectives on the nature of reality, consciousness, and the relationship between mind and matter" " By delving into the past, we gain invaluable insights into the present, fostering a dialogue between the two that enriches our understanding of ourselves and the world around us"

# Remove the final sentence + trailing period of the summary paragraph entirely
Replace-Text " While quantum theory's complexities continue to intrigue and challenge our intellect, they also hold the promise of unlocking hidden realms of scientific knowledge and technological advancements, pushing the boundaries of human understanding and innovation." ""

# --- 5. Append a new empty paragraph at the end of the document ---
$end = $d.Range($d.Content.End, $d.Content.End)
$end.InsertParagraphAfter()

Write-Output "done"
